$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text corrections (shared string content) ---
$ws.Range("C1").Value = "Name"
$ws.Range("I1").Value = "Booking#"
$ws.Range("J1").Value = "Registration#"
$ws.Range("K1").Value = "Registration Type"

# --- Row heights: 16.5 -> 18.75 for header + data row ---
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(2).RowHeight = 18.75

# --- Column K width: widen to fit the new "Registration#" header ---
$ws.Columns.Item(11).ColumnWidth = 16.83

# --- Font color: theme color -> explicit black (rgb FF000000) on the
#     cells that used the bordered/header font (fontId 1), i.e. every
#     styled data cell except the 4 blank "general" cells in row 2. ---
$ws.Range("A1:D1").Font.Color = 0
$ws.Range("E1:X1").Font.Color = 0
$ws.Range("Y1:AB1").Font.Color = 0
$ws.Range("A2:D2").Font.Color = 0
$ws.Range("G2:X2").Font.Color = 0
$ws.Range("Z2:AA2").Font.Color = 0

# --- Leave a spare "general left-aligned" style allocated (matches the
#     extra cellXfs entry introduced upstream) without leaving it
#     attached to any cell. ---
$ws.Range("E2").HorizontalAlignment = -4131
$ws.Range("E2").HorizontalAlignment = 1

Write-Output "done"
